{"js": "// Add a new comment (\"Yes, it matches, and is posted correctly in Canvas.\n// Thank you.\") anchored to the word \"Memorandum\" at the very top of the\n// memo -- the same word that already carries the original review\n// comment from Cannell, Michael B.\nconst body = context.document.body;\nconst results = body.search(\"Memorandum\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the text 'Memorandum' to attach the comment to.\");\n}\n\n// The heading \"Memorandum\" occurs once, at the top of the document --\n// grab that occurrence and insert the new comment on it.\nconst target = results.items[0];\ntarget.insertComment(\"Yes, it matches, and is posted correctly in Canvas. Thank you.\");\n\nawait context.sync();\n", "ps1": "# Add a new review comment -- \"Yes, it matches, and is posted correctly in\n# Canvas. Thank you.\" -- anchored to the word \"Memorandum\" at the top of\n# the memo, authored by Callender, Librada (the same reviewer already\n# replying in-line in the existing comment thread).\n$d = $word.ActiveDocument\n\n# The new comment is from Librada Callender, so make sure the identity\n# Word stamps onto the new comment matches her (not whatever the COM\n# host's default identity happens to be).\n$word.UserName = \"Callender, Librada\"\n$word.UserInitials = \"CL\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Memorandum\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find the text 'Memorandum' to attach the comment to.\"\n}\n\n$d.Comments.Add($find.Parent, \"Yes, it matches, and is posted correctly in Canvas. Thank you.\") | Out-Null\n"}
